$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update comment text in E4 and E5 to include extra detail
$ws.Range("E4").Value = "Consistency in results for both winning teams and losing teams, Class 0 and Class 1 have 614"
$ws.Range("E5").Value = "Optimized Model: Consistency in results for both winning teams and losing teams, Class 0 and Class 1 have 614"

# Widen column E to fit the longer comments
$ws.Range("E1:E1").EntireColumn.ColumnWidth = 99

# Update the active selection to E6
$ws.Range("E6").Select()
